# 螺纹钢销量.xlsx — append two new monthly data refreshes (2023-12 / 2024-01)
# pulled in by the EM_EDB_N (Choice) data feed, matching what a live refresh
# of the linked EDB range would have written below the prior last row, and
# re-append the trailing "data source" footer note one row further down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 165: 2023-12-31 refresh -------------------------------------------------
# (date written as its serial number via Value2 so Excel doesn't silently mint
# a brand-new "m/d/yy" number-format style from the string parse — the cell
# picks up the workbook's existing yyyy-mm-dd style via the format copy below)
$ws.Rows(165).RowHeight = 14
$ws.Range("A165").Value2 = 45291
$ws.Range("B165").Value = 49312586
$ws.Range("C165").Value = 8367.1831000000002
$ws.Range("D165").Value = 161784191
$ws.Range("E165").Value = 8564342
$ws.Range("F165").Value = 1938.7186999999999
$ws.Range("G165").Value = 84.824299999999994

# --- Row 166: 2024-01-31 refresh -------------------------------------------------
$ws.Rows(166).RowHeight = 14
$ws.Range("A166").Value2 = 45322
$ws.Range("B166").Value = 3194666
$ws.Range("C166").Value = 603.88549999999998
$ws.Range("D166").Value = 11104838
$ws.Range("E166").Value = 552561
$ws.Range("F166").Value = 124.8873
$ws.Range("G166").Value = 6.9882

# Give the two new data rows the same look as the data rows above them
# (date format on col A, numeric format on cols B:G).
$ws.Range("A164").Copy()
$ws.Range("A165:A166").PasteSpecial(-4122)
$ws.Range("B164:G164").Copy()
$ws.Range("B165:G166").PasteSpecial(-4122)

# --- Row 169: re-stamp the trailing source-note row below the new data ----------
# (row 167 already carries it; row 168 is intentionally left blank, matching
# the gap produced by the EDB refresh pushing the note down twice)
$ws.Range("A167").Copy()
$ws.Range("A169").PasteSpecial(-4122)
$ws.Range("A169").Value = $ws.Range("A167").Value2

$excel.CutCopyMode = 0

# Selection left where the refresh's last interactive edit was made.
$null = $ws.Range("C163").Select()
